# Fix outfall bug: move the "Tidal" curve sheet so it sits right after
# "Storage" (ahead of Weir/Pump1/Pump3/Pump4), and populate it with the
# tidal-stage curve data that was previously missing (empty sheet).

$wb = $excel.ActiveWorkbook

# --- Reorder: move "Tidal" tab to just before "Weir" -----------------
$tidal = $wb.Worksheets.Item("Tidal")
$weir  = $wb.Worksheets.Item("Weir")
$tidal.Move($weir)

# --- Populate the Tidal curve data ------------------------------------
$ws = $wb.Worksheets.Item("Tidal")

$data = @(
    @("tdc_1", 1, 0),
    @("tdc_1", 2, 0),
    @("tdc_1", 3, 0.3),
    @("tdc_1", 4, 0.2),
    @("tdc_1", 5, 0.2),
    @("tdc_1", 6, 0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = ""
    $r = $r + 1
}
